$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 358
$lastNewRow = 366
$lastOldRow = 357

# Copy the formatting of the last existing row down onto the new rows so the
# same style indices (e.g. date format on column A) are reused instead of
# new styles being created.
$ws.Cells.Item($lastOldRow, 1).Copy() | Out-Null
$ws.Range($ws.Cells.Item($firstNewRow, 1), $ws.Cells.Item($lastNewRow, 1)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data: serial date, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44432, 0, 0, 0),
    @(44433, 0, 0, 0),
    @(44434, 0, 0, 0),
    @(44435, 0, 0, 0),
    @(44436, 0, 0, 0),
    @(44437, 0, 0, 0),
    @(44438, 2, 2, 61.06870229007634),
    @(44439, 0, 2, 61.06870229007634),
    @(44440, 0, 2, 61.06870229007634)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $firstNewRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
